# Auto-generated Excel COM-interop script
# Applies numeric value updates (Tonberry_Profits market-data refresh) across 8 sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 780.61536
$ws.Range("I15").Value = 780.61536
$ws.Range("K15").Value = 2341.84608
$ws.Range("M15").Value = -2172.84608
# row 28
$ws.Range("H28").Value = 1196
$ws.Range("I28").Value = 145.33333
$ws.Range("J28").Value = 7500
$ws.Range("K28").Value = 145.33333
$ws.Range("L28").Value = 7500
$ws.Range("M28").Value = 339.66667
$ws.Range("N28").Value = -8470
# row 62
$ws.Range("H62").Value = 2734.6667
$ws.Range("I62").Value = 2499.6667
$ws.Range("J62").Value = 2969.6667
$ws.Range("K62").Value = 2499.6667
$ws.Range("L62").Value = 2969.6667
$ws.Range("M62").Value = -1875.6667
$ws.Range("N62").Value = -4217.6667
# row 65
$ws.Range("H65").Value = 2734.6667
$ws.Range("I65").Value = 2499.6667
$ws.Range("J65").Value = 2969.6667
$ws.Range("K65").Value = 12498.3335
$ws.Range("L65").Value = 14848.3335
$ws.Range("M65").Value = -9378.333500000001
$ws.Range("N65").Value = -21088.3335
# row 75
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31872
# row 78
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99360
# row 98
$ws.Range("H98").Value = 2431.3
$ws.Range("J98").Value = 3399
$ws.Range("L98").Value = 3399
$ws.Range("N98").Value = -6395
# row 121
$ws.Range("H121").Value = 1500
$ws.Range("J121").Value = 1500
$ws.Range("L121").Value = 4500
$ws.Range("N121").Value = -7994
# row 122
$ws.Range("H122").Value = 2431.3
$ws.Range("J122").Value = 3399
$ws.Range("L122").Value = 10197
$ws.Range("N122").Value = -15097
# row 132
$ws.Range("H132").Value = 1068.0264
$ws.Range("I132").Value = 1015.8108
$ws.Range("K132").Value = 3047.4324
$ws.Range("M132").Value = -517.4323999999997
# row 137
$ws.Range("H137").Value = 1818.3572
$ws.Range("J137").Value = 1935.5
$ws.Range("L137").Value = 5806.5
$ws.Range("N137").Value = -10906.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 3602.0508
$ws.Range("I32").Value = 2223.8223
$ws.Range("K32").Value = 2223.8223
$ws.Range("M32").Value = -1936.8223
# row 45
$ws.Range("H45").Value = 1543.1428
$ws.Range("I45").Value = 925.6875
$ws.Range("K45").Value = 925.6875
$ws.Range("M45").Value = -548.6875
# row 61
$ws.Range("H61").Value = 4822.1665
$ws.Range("I61").Value = 2196
$ws.Range("K61").Value = 2196
$ws.Range("M61").Value = -1984
# row 74
$ws.Range("H74").Value = 10000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
# row 77
$ws.Range("H77").Value = 10000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
# row 110
$ws.Range("H110").Value = 1685.2727
$ws.Range("I110").Value = 1616.4445
$ws.Range("J110").Value = 1995
$ws.Range("K110").Value = 1616.4445
$ws.Range("L110").Value = 1995
$ws.Range("M110").Value = 428.5554999999999
$ws.Range("N110").Value = -6085
# row 132
$ws.Range("H132").Value = 1986.5
$ws.Range("I132").Value = 1428.0555
$ws.Range("K132").Value = 4284.166499999999
$ws.Range("M132").Value = -1754.166499999999
# row 136
$ws.Range("H136").Value = 4822.1665
$ws.Range("I136").Value = 2196
$ws.Range("K136").Value = 6588
$ws.Range("M136").Value = -4038

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 105
$ws.Range("H105").Value = 2399.8125
$ws.Range("I105").Value = 2399.8125
$ws.Range("K105").Value = 2399.8125
$ws.Range("M105").Value = -652.8125
# row 134
$ws.Range("H134").Value = 8556.923000000001
$ws.Range("I134").Value = 9336.65
$ws.Range("J134").Value = 5957.8335
$ws.Range("K134").Value = 28009.95
$ws.Range("L134").Value = 17873.5005
$ws.Range("M134").Value = -25474.95
$ws.Range("N134").Value = -22943.5005

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 2000
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2700
# row 31
$ws.Range("H31").Value = 3954.6316
$ws.Range("I31").Value = 1440.25
$ws.Range("K31").Value = 1440.25
$ws.Range("M31").Value = -1145.25
# row 34
$ws.Range("H34").Value = 3954.6316
$ws.Range("I34").Value = 1440.25
$ws.Range("K34").Value = 1440.25
$ws.Range("M34").Value = -1238.25
# row 63
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# row 66
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# row 132
$ws.Range("H132").Value = 2177.68
$ws.Range("I132").Value = 1120.8
$ws.Range("K132").Value = 3362.4
$ws.Range("M132").Value = -832.3999999999996
# row 134
$ws.Range("H134").Value = 1092.4166
$ws.Range("I134").Value = 1059.6
$ws.Range("K134").Value = 3178.8
$ws.Range("M134").Value = -643.7999999999997

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 107
$ws.Range("H107").Value = 560.94446
$ws.Range("J107").Value = 560.94446
$ws.Range("L107").Value = 1682.83338
$ws.Range("N107").Value = -5522.83338
# row 131
$ws.Range("H131").Value = 11924452
$ws.Range("J131").Value = 24193.766
$ws.Range("L131").Value = 72581.298
$ws.Range("N131").Value = -82661.298

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 1567.2858
$ws.Range("I102").Value = 1371.625
$ws.Range("K102").Value = 1371.625
$ws.Range("M102").Value = 250.375
# row 132
$ws.Range("H132").Value = 2756.1082
$ws.Range("I132").Value = 2437.6453
$ws.Range("K132").Value = 7312.9359
$ws.Range("M132").Value = -4782.9359
# row 134
$ws.Range("H134").Value = 46999.4
$ws.Range("J134").Value = 46999.4
$ws.Range("L134").Value = 140998.2
$ws.Range("N134").Value = -146068.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 1450.6666
$ws.Range("J22").Value = 1463.4166
$ws.Range("L22").Value = 1463.4166
$ws.Range("N22").Value = -2053.4166
# row 27
$ws.Range("H27").Value = 1450.6666
$ws.Range("J27").Value = 1463.4166
$ws.Range("L27").Value = 1463.4166
$ws.Range("N27").Value = -1677.4166
# row 46
$ws.Range("H46").Value = 1936.3636
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1936.3636
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1936.3636
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -2312.3636
# row 68
$ws.Range("H68").Value = 3362.1428
$ws.Range("J68").Value = 3777.5
$ws.Range("L68").Value = 3777.5
$ws.Range("N68").Value = -5275.5
# row 71
$ws.Range("H71").Value = 3362.1428
$ws.Range("J71").Value = 3777.5
$ws.Range("L71").Value = 18887.5
$ws.Range("N71").Value = -26375.5
# row 132
$ws.Range("H132").Value = 1777.0646
$ws.Range("I132").Value = 1364.8462
$ws.Range("J132").Value = 2074.7778
$ws.Range("K132").Value = 4094.5386
$ws.Range("L132").Value = 6224.3334
$ws.Range("M132").Value = -1564.5386
$ws.Range("N132").Value = -11284.3334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376
# row 65
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880
# row 81
$ws.Range("H81").Value = 500
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 500
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 1000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -3122
# row 84
$ws.Range("H84").Value = 500
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 500
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 5000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -15608
# row 113
$ws.Range("H113").Value = 828
# row 122
$ws.Range("H122").Value = 15774.172
$ws.Range("I122").Value = 21345.36
$ws.Range("K122").Value = 64036.08
$ws.Range("M122").Value = -61586.08
